$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '25.199.37'
$ws.Range("E2").Value = '  -2.97%  '

# Row 3
$ws.Range("D3").Value = '1.549.02'
$ws.Range("E3").Value = '  -4.79%  '

# Row 4
$ws.Range("E4").Value = '  +0.03%  '

# Row 5
$ws.Range("D5").Value = '206.34'
$ws.Range("E5").Value = '  -3.53%  '

# Row 6
$ws.Range("E6").Value = '  +0.24%  '

# Row 7
$ws.Range("D7").Value = '0.475'
$ws.Range("E7").Value = '  -5.59%  '

# Row 8
$ws.Range("D8").Value = '0.0604'
$ws.Range("E8").Value = '  -2.21%  '

# Row 9
$ws.Range("D9").Value = '0.240'
$ws.Range("E9").Value = '  -3.81%  '

# Row 10
$ws.Range("D10").Value = '17.69'
$ws.Range("E10").Value = '  -3.62%  '

# Row 11
$ws.Range("D11").Value = '0.0780'
$ws.Range("E11").Value = '  -1.02%  '

# Row 12
$ws.Range("D12").Value = '1.764.22'
$ws.Range("E12").Value = '  -4.64%  '

# Row 13
$ws.Range("D13").Value = '1.553.04'
$ws.Range("E13").Value = '  -5.95%  '

# Row 14
$ws.Range("D14").Value = '3.96'
$ws.Range("E14").Value = '  -5.28%  '

# Row 15
$ws.Range("D15").Value = '0.500'
$ws.Range("E15").Value = '  -4.82%  '

# Row 16
$ws.Range("D16").Value = '25.133.66'
$ws.Range("E16").Value = '  -3.05%  '

# Row 17
$ws.Range("D17").Value = '0.0₃0707'
$ws.Range("E17").Value = '  -4.02%  '

# Row 18
$ws.Range("D18").Value = '58.48'
$ws.Range("E18").Value = '  -4.50%  '

# Row 19
$ws.Range("E19").Value = '  -0.01%  '

# Row 20
$ws.Range("D20").Value = '184.54'
$ws.Range("E20").Value = '  -3.81%  '

# Row 21
$ws.Range("E21").Value = '  -3.61%  '

# Row 22
$ws.Range("D22").Value = '9.19'
$ws.Range("E22").Value = '  -4.06%  '

# Row 23
$ws.Range("D23").Value = '5.81'
$ws.Range("E23").Value = '  -4.24%  '

# Row 24
$ws.Range("E24").Value = '  -0.13%  '

# Row 25
$ws.Range("B25").Value = 'Stellar'
$ws.Range("C25").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D25").Value = '0.127'
$ws.Range("E25").Value = '  -4.46%  '

# Row 26
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").Value = '139.08'
$ws.Range("E26").Value = '  -3.14%  '

# Row 27
$ws.Range("D27").Value = '1.65'
$ws.Range("E27").Value = '  -4.66%  '

# Row 28
$ws.Range("D28").Value = '14.74'
$ws.Range("E28").Value = '  -2.99%  '

# Row 29
$ws.Range("E29").Value = '  -5.13%  '

# Row 30
$ws.Range("D30").Value = '1.14'
$ws.Range("E30").Value = '  -6.94%  '

# Row 31
$ws.Range("D31").Value = '0.0460'
$ws.Range("E31").Value = '  -4.82%  '

# Row 32
$ws.Range("D32").Value = '3.01'
$ws.Range("E32").Value = '  -3.74%  '

# Row 33
$ws.Range("D33").Value = '2.96'
$ws.Range("E33").Value = '  -4.58%  '

# Row 34
$ws.Range("E34").Value = '  -3.98%  '

# Row 35
$ws.Range("E35").Value = '  -3.90%  '

# Row 36
$ws.Range("D36").Value = '1.085.24'
$ws.Range("E36").Value = '  -3.25%  '

# Row 37
$ws.Range("E37").Value = '  -0.36%  '

# Row 38
$ws.Range("E38").Value = '  -2.75%  '

# Row 39
$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").Value = '0.489'
$ws.Range("E39").Value = '  -5.74%  '

# Row 40
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").Value = '2.25'
$ws.Range("E40").Value = '  -7.46%  '

# Row 41
$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D41").Value = '0.756'
$ws.Range("E41").Value = '  -10.92%  '

# Row 42
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = '0.793'
$ws.Range("E42").Value = '  +3.55%  '

# Row 43
$ws.Range("D43").Value = '92.12'
$ws.Range("E43").Value = '  -5.99%  '

# Row 44
$ws.Range("D44").Value = '5.01'
$ws.Range("E44").Value = '  -2.89%  '

# Row 45
$ws.Range("D45").Value = '1.679.70'
$ws.Range("E45").Value = '  -4.58%  '

# Row 46
$ws.Range("D46").Value = '0.0₆0106'
$ws.Range("E46").Value = '  +1.42%  '

# Row 47
$ws.Range("D47").Value = '52.09'
$ws.Range("E47").Value = '  -4.10%  '

# Row 48
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").Value = '1.44'
$ws.Range("E48").Value = '  -2.16%  '

# Row 49
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '0.0501'
$ws.Range("E49").Value = '  -5.58%  '

# Row 50
$ws.Range("E50").Value = '  -1.77%  '

# Row 51
$ws.Range("E51").Value = '  -0.32%  '
